$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(0,"falling",-1.765928411483764,1.167147111892701,0.571823143959045,-0.0310014113783836,0.0154243474826216,0.0794124826788902),
    @(100,"falling",-1.709035491943359,1.218545722961426,0.4037320613861088,-0.005192354787141,0.0403171069920063,-0.0171042270958423),
    @(200,"falling",-1.738995742797851,1.102185392379761,0.6806826651096344,0.0421497002243995,-0.0251981914043426,0.0042760567739605),
    @(300,"falling",-1.775501251220703,1.23273515701294,0.5835052132606506,0.0103847095742821,-0.0126754539087414,-0.00335975876078),
    @(400,"falling",-1.910298776626587,1.120057487487793,0.7890581786632539,0.0180205255746841,-0.0375682115554809,0.0152716310694813),
    @(500,"falling",-1.610696506500243,1.25785665512085,0.8155304193496704,-0.0210748501121997,-0.0261144898831844,-0.0280998013913631),
    @(600,"falling",-2.127863931655885,1.3188805103302,0.8300568103790283,-0.0584903471171855,0.0059559359215199,0.0522289797663688),
    @(700,"falling",-2.092850303649902,0.9900987625122062,1.828120517730715,0.0937678143382072,0.0774271711707115,-0.1965458989143371),
    @(800,"falling",-2.105730056762695,0.8735208511352541,2.453470587730407,-0.1837177276611328,0.2562579810619354,0.0125227374956011),
    @(900,"falling",-2.657281875610353,0.2476238250732379,2.285648679733275,-0.0036651915870606,0.1434006094932556,-0.0604756586253643),
    @(1000,"falling",-4.063823699951181,-1.187473011016849,2.541466045379643,0.2081523388624191,0.1316414624452591,0.1007927656173706),
    @(1100,"falling",-5.616066837310788,-2.403632545471194,3.652213478088379,-0.4978551864624023,-0.2622139155864715,-1.034500241279602),
    @(1200,"falling",-3.539181804656983,-3.395954704284656,3.008686828613268,-0.328340083360672,-0.4218024611473083,-0.5737552046775818),
    @(1300,"falling",-3.633686542510985,3.542618751525879,-4.68873119354248,-2.957962274551392,-0.2176207453012466,-2.633592844009399),
    @(1400,"falling",53.59859118461612,-4.54958009719849,-34.16661596298219,1.301142930984497,-3.674507141113281,2.831360340118408),
    @(1500,"falling",10.96323089599605,-0.1360907554626412,-11.01923332214353,0.7254024744033813,0.9521862268447876,-1.527926683425903),
    @(1600,"falling",-1.288498878479004,0.4897777378559073,-1.382943773269636,-0.2622139155864715,1.557553648948669,-0.6866125464439392),
    @(1700,"falling",0.2403127670288223,0.4333343148231529,1.078743743896495,1.219134330749511,-0.9372199773788452,2.843730449676514),
    @(1800,"falling",0.6345248222351046,0.4263583719730371,1.455561161041258,0.8478809595108032,-0.9918924570083618,-0.8237518072128296),
    @(1900,"falling",-0.09080390930175844,2.436175584793104,1.085069131851197,-0.0302378293126821,-0.6291912198066711,-0.0180205255746841),
    @(2000,"falling",-1.009081411361702,3.721519541740399,1.14239126443863,-0.0595593601465225,-0.2151772826910019,0.1901318132877349),
    @(2100,"falling",-2.444419670104975,-0.2522135257720842,1.136495351791381,0.0241291765123605,0.3474296033382416,0.1518000066280365),
    @(2200,"falling",-1.008091449737549,2.286790394783017,0.9021725535392736,-0.0592539273202419,0.4436408877372741,0.2559525370597839),
    @(2300,"falling",-1.014060974121094,1.50408124923706,0.3019133806228637,0.057115901261568,0.0777326002717018,-0.1259909570217132),
    @(2400,"falling",-0.6290699958801254,1.975620031356814,-0.5599067687988319,-0.0335975885391235,-0.0216857157647609,-0.110566608607769),
    @(2500,"falling",-1.261460304260257,1.682102203369138,-0.007327961921688053,0.0045814891345798,0.0555887371301651,0.012980886735022),
    @(2600,"falling",-0.9165861129760703,1.882461953163149,0.1453821629285811,-0.0044287731871008,0.0375682115554809,0.0296269636601209),
    @(2700,"falling",-0.8792205810546895,1.9060809135437,0.1410199522972106,-0.0070249503478407,-0.0122173046693205,0.0006108652451075),
    @(2800,"falling",-0.9211750030517577,1.817085385322571,0.09843596816062908,-0.0004581489483825,0.0073303831741213,0.0157297793775796),
    @(2900,"falling",-0.8201595306396483,1.820924615859986,0.07434962689876565,-0.0102319931611418,-0.020616702735424,0.0047342055477201)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
